$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.816.15"
$ws.Range("D3").Value = "1.633.35"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "'214.76"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'0.5019"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "'0.2565"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").Value = "'0.06391"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").Value = "'19.60"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("D11").Value = "'0.07693"
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.241"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.633.78"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "1.858.25"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "'0.5415"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").Value = "0.0₅7906"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "'63.47"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "25.839.88"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "'202.02"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("D21").Value = "'4.326"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "'9.916"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").Value = "'5.944"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "'1.933"
$ws.Range("E25").Value = "  +11.90%  "
$ws.Range("D26").Value = "'141.46"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").Value = "'0.1140"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").Value = "'15.66"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  -3.83%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "'0.04988"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").Value = "'3.256"
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("D36").Value = "1.170.84"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").Value = "'0.8903"
$ws.Range("E37").Value = "  -3.99%  "
$ws.Range("D38").Value = "'2.614"
$ws.Range("E38").Value = "  -4.93%  "
$ws.Range("D39").Value = "'0.5573"
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").Value = "'0.8066"
$ws.Range("E44").Value = "  -3.26%  "
$ws.Range("D45").Value = "'99.24"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "1.770.41"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").Value = "'0.4513"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "'54.61"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("E51").Value = "  +0.75%  "
